$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 9, shifting the existing rows 9-24 down to 10-25
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the new weekly record (same constant columns as
# the surrounding rows, new date/volume/price figures)
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 45079
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112001
$ws.Range("G9").Value = "Berenjena"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 130
$ws.Range("K9").Value = 4000
$ws.Range("L9").Value = 5000
$ws.Range("M9").Value = 4462
$ws.Range("N9").Value = "$/caja 60 unidades"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 74
$ws.Range("Q9").Value = 60
$ws.Range("R9").Value = "Hortaliza"
